$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Daily Orders" - a brand-new order (Order ID 12) comes in, so it is
#    inserted as the new row 2 (most-recent-first ordering), pushing all
#    existing order rows down by one.
# ---------------------------------------------------------------------------
$orders = $wb.Worksheets.Item("Daily Orders")
$orders.Rows.Item(2).Insert()

# D (Flat No) and J (Collection Date) hold digit-only / date-shaped text in
# this sheet ("420", "2026-01-15"); force Text format first so Excel doesn't
# silently reinterpret them as a number / date serial.
$orders.Range("D2").NumberFormat = "@"
$orders.Range("J2").NumberFormat = "@"

$orders.Range("A2").Value2 = 12
$orders.Range("B2").Value2 = "2026-01-13 22:43"
$orders.Range("C2").Value2 = "Swapnil (Phantom)"
$orders.Range("D2").Value2 = "420"
$orders.Range("E2").Value2 = ""
$orders.Range("F2").Value2 = "Vermicelli Kheer x1"
$orders.Range("G2").Value2 = 50
$orders.Range("H2").Value2 = "NEW"
$orders.Range("I2").Value2 = "PENDING"
$orders.Range("J2").Value2 = "2026-01-15"
$orders.Range("K2").Value2 = "16:42"
$orders.Range("L2").Value2 = "No vermicelli in kheer please."
$orders.Range("M2").Value2 = ""
$orders.Range("N2").Value2 = ""

# Drop back to the sheet's normal (unformatted) style now that the text is
# locked in, so D2/J2 don't carry a lingering explicit "Text" number format.
$orders.Range("D2").Style = "Normal"
$orders.Range("J2").Style = "Normal"

# ---------------------------------------------------------------------------
# 2) "Summary" - totals reflect the new order: one more total order, one more
#    "New" order, and the revenue total grows by the order's 50.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("A2").Value2 = 12
$summary.Range("B2").Value2 = 10
$summary.Range("G2").Value2 = 325

# ---------------------------------------------------------------------------
# 3) "Items Breakdown" - add the new item "Vermicelli Kheer" in its sorted
#    position (row 4, right after "Til Poli"), pushing the remaining items
#    down by one row.
# ---------------------------------------------------------------------------
$breakdown = $wb.Worksheets.Item("Items Breakdown")
$breakdown.Rows.Item(4).Insert()

$breakdown.Range("A4").Value2 = "Vermicelli Kheer"
$breakdown.Range("B4").Value2 = 1
$breakdown.Range("C4").Value2 = 50
